$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 39945
$ws.Range("J3").Value = 39945
$ws.Range("L3").Value = 39945
$ws.Range("N3").Value = -40173
$ws.Range("H5").Value = 28.5
$ws.Range("J5").Value = 8
$ws.Range("L5").Value = 8
$ws.Range("N5").Value = -238
$ws.Range("H43").Value = 4937.3335
$ws.Range("I43").Value = 2311
$ws.Range("J43").Value = 6250.5
$ws.Range("K43").Value = 2311
$ws.Range("L43").Value = 6250.5
$ws.Range("M43").Value = -2242
$ws.Range("N43").Value = -6388.5
$ws.Range("H62").Value = 19233500
$ws.Range("I62").Value = 22730096
$ws.Range("K62").Value = 22730096
$ws.Range("M62").Value = -22729472
$ws.Range("H65").Value = 19233500
$ws.Range("I65").Value = 22730096
$ws.Range("K65").Value = 113650480
$ws.Range("M65").Value = -113647360
$ws.Range("H96").Value = 802.4
$ws.Range("I96").Value = 689.5
$ws.Range("J96").Value = 877.6667
$ws.Range("K96").Value = 2068.5
$ws.Range("L96").Value = 2633.0001
$ws.Range("M96").Value = -695.5
$ws.Range("N96").Value = -5379.0001
$ws.Range("H98").Value = 606.2857
$ws.Range("I98").Value = 587
$ws.Range("J98").Value = 722
$ws.Range("K98").Value = 587
$ws.Range("L98").Value = 722
$ws.Range("M98").Value = 911
$ws.Range("N98").Value = -3718
$ws.Range("H100").Value = 1698.7
$ws.Range("I100").Value = 1453.1666
$ws.Range("J100").Value = 2067
$ws.Range("K100").Value = 1453.1666
$ws.Range("L100").Value = 2067
$ws.Range("M100").Value = -912.1666
$ws.Range("N100").Value = -3149
$ws.Range("H102").Value = 39945
$ws.Range("J102").Value = 39945
$ws.Range("L102").Value = 39945
$ws.Range("N102").Value = -46435
$ws.Range("H103").Value = 1863.8572
$ws.Range("J103").Value = 449.33334
$ws.Range("L103").Value = 1348.00002
$ws.Range("N103").Value = -2520.00002
$ws.Range("H113").Value = 5787.391
$ws.Range("I113").Value = 3959.5
$ws.Range("J113").Value = 6762.2666
$ws.Range("K113").Value = 3959.5
$ws.Range("L113").Value = 6762.2666
$ws.Range("M113").Value = -705.5
$ws.Range("N113").Value = -13270.2666
$ws.Range("H116").Value = 11667.9
$ws.Range("I116").Value = 3394.6
$ws.Range("J116").Value = 19941.2
$ws.Range("K116").Value = 3394.6
$ws.Range("L116").Value = 19941.2
$ws.Range("M116").Value = 47.40000000000009
$ws.Range("N116").Value = -26825.2
$ws.Range("H122").Value = 606.2857
$ws.Range("I122").Value = 587
$ws.Range("J122").Value = 722
$ws.Range("K122").Value = 1761
$ws.Range("L122").Value = 2166
$ws.Range("M122").Value = 689
$ws.Range("N122").Value = -7066
$ws.Range("H129").Value = 3412.875
$ws.Range("I129").Value = 1100.75
$ws.Range("J129").Value = 5725
$ws.Range("K129").Value = 3302.25
$ws.Range("L129").Value = 17175
$ws.Range("M129").Value = 1697.75
$ws.Range("N129").Value = -27175
$ws.Range("H132").Value = 3357.561
$ws.Range("I132").Value = 3591.862
$ws.Range("J132").Value = 2791.3333
$ws.Range("K132").Value = 10775.586
$ws.Range("L132").Value = 8373.999899999999
$ws.Range("M132").Value = -8245.585999999999
$ws.Range("N132").Value = -13433.9999
$ws.Range("H135").Value = 7682.6445
$ws.Range("I135").Value = 838.8205
$ws.Range("J135").Value = 52167.5
$ws.Range("K135").Value = 7549.3845
$ws.Range("L135").Value = 469507.5
$ws.Range("M135").Value = -5014.3845
$ws.Range("N135").Value = -474577.5
$ws.Range("H137").Value = 4538.5
$ws.Range("I137").Value = 3877.6
$ws.Range("J137").Value = 5199.4
$ws.Range("K137").Value = 11632.8
$ws.Range("L137").Value = 15598.2
$ws.Range("M137").Value = -9082.799999999999
$ws.Range("N137").Value = -20698.2
$ws.Range("H141").Value = 5882.8423
$ws.Range("I141").Value = 4173.4375
$ws.Range("J141").Value = 14999.667
$ws.Range("K141").Value = 12520.3125
$ws.Range("L141").Value = 44999.001
$ws.Range("M141").Value = -7340.3125
$ws.Range("N141").Value = -55359.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5436.2173
$ws.Range("I32").Value = 5096.143
$ws.Range("K32").Value = 5096.143
$ws.Range("M32").Value = -4809.143
$ws.Range("H61").Value = 70007
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H74").Value = 1951.3334
$ws.Range("I74").Value = 1951.3334
$ws.Range("K74").Value = 1951.3334
$ws.Range("M74").Value = -1077.3334
$ws.Range("H77").Value = 1951.3334
$ws.Range("I77").Value = 1951.3334
$ws.Range("K77").Value = 9756.666999999999
$ws.Range("M77").Value = -5388.666999999999
$ws.Range("H88").Value = 35740.5
$ws.Range("I88").Value = 25396.334
$ws.Range("J88").Value = 46084.668
$ws.Range("K88").Value = 25396.334
$ws.Range("L88").Value = 46084.668
$ws.Range("M88").Value = -24990.334
$ws.Range("N88").Value = -46896.668
$ws.Range("H91").Value = 35740.5
$ws.Range("I91").Value = 25396.334
$ws.Range("J91").Value = 46084.668
$ws.Range("K91").Value = 25396.334
$ws.Range("L91").Value = 46084.668
$ws.Range("M91").Value = -23992.334
$ws.Range("N91").Value = -48892.668
$ws.Range("I102").Value = 1044.3846
$ws.Range("J102").Value = 1000
$ws.Range("K102").Value = 1044.3846
$ws.Range("L102").Value = 1000
$ws.Range("M102").Value = 577.6153999999999
$ws.Range("N102").Value = -4244
$ws.Range("H110").Value = 4173.68
$ws.Range("I110").Value = 4238.5454
$ws.Range("K110").Value = 4238.5454
$ws.Range("M110").Value = -2193.5454
$ws.Range("H122").Value = 1900
$ws.Range("I122").Value = 1850
$ws.Range("K122").Value = 5550
$ws.Range("M122").Value = -3100
$ws.Range("H132").Value = 3642.5264
$ws.Range("I132").Value = 3606.4119
$ws.Range("K132").Value = 10819.2357
$ws.Range("M132").Value = -8289.235700000001
$ws.Range("H136").Value = 70007
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2697.25
$ws.Range("I20").Value = 1395
$ws.Range("K20").Value = 1395
$ws.Range("M20").Value = -1148
$ws.Range("H22").Value = 3115
$ws.Range("I22").Value = 3115
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 3115
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -2942
$ws.Range("H86").Value = 24246.95
$ws.Range("I86").Value = 6224.6665
$ws.Range("J86").Value = 51280.375
$ws.Range("K86").Value = 6224.6665
$ws.Range("L86").Value = 51280.375
$ws.Range("M86").Value = -5101.6665
$ws.Range("N86").Value = -53526.375
$ws.Range("H89").Value = 24246.95
$ws.Range("I89").Value = 6224.6665
$ws.Range("J89").Value = 51280.375
$ws.Range("K89").Value = 31123.3325
$ws.Range("L89").Value = 256401.875
$ws.Range("M89").Value = -25507.3325
$ws.Range("N89").Value = -267633.875
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("H99").Value = 2254.889
$ws.Range("I99").Value = 2254.889
$ws.Range("K99").Value = 2254.889
$ws.Range("M99").Value = -756.8890000000001
$ws.Range("H107").Value = 2421.2173
$ws.Range("I107").Value = 2224.1
$ws.Range("J107").Value = 3735.3333
$ws.Range("K107").Value = 2224.1
$ws.Range("L107").Value = 3735.3333
$ws.Range("M107").Value = -304.0999999999999
$ws.Range("N107").Value = -7575.3333
$ws.Range("H134").Value = 6536.6445
$ws.Range("I134").Value = 6248.1353
$ws.Range("J134").Value = 7871
$ws.Range("K134").Value = 18744.4059
$ws.Range("L134").Value = 23613
$ws.Range("M134").Value = -16209.4059
$ws.Range("N134").Value = -28683

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 699.2
$ws.Range("I16").Value = 573.4286
$ws.Range("J16").Value = 992.6667
$ws.Range("K16").Value = 573.4286
$ws.Range("L16").Value = 992.6667
$ws.Range("M16").Value = -286.4286
$ws.Range("N16").Value = -1566.6667
$ws.Range("H22").Value = 963.1818
$ws.Range("I22").Value = 713.5714
$ws.Range("J22").Value = 1400
$ws.Range("K22").Value = 713.5714
$ws.Range("L22").Value = 1400
$ws.Range("M22").Value = -363.5714
$ws.Range("N22").Value = -2100
$ws.Range("H31").Value = 3281.7046
$ws.Range("I31").Value = 2842.3333
$ws.Range("J31").Value = 3394.6858
$ws.Range("K31").Value = 2842.3333
$ws.Range("L31").Value = 3394.6858
$ws.Range("M31").Value = -2547.3333
$ws.Range("N31").Value = -3984.6858
$ws.Range("H34").Value = 3281.7046
$ws.Range("I34").Value = 2842.3333
$ws.Range("J34").Value = 3394.6858
$ws.Range("K34").Value = 2842.3333
$ws.Range("L34").Value = 3394.6858
$ws.Range("M34").Value = -2640.3333
$ws.Range("N34").Value = -3798.6858
$ws.Range("H58").Value = 3522.077
$ws.Range("I58").Value = 2871.5908
$ws.Range("K58").Value = 2871.5908
$ws.Range("M58").Value = -2668.5908
$ws.Range("H62").Value = 22025.166
$ws.Range("I62").Value = 5777.5
$ws.Range("K62").Value = 5777.5
$ws.Range("M62").Value = -5153.5
$ws.Range("H65").Value = 22025.166
$ws.Range("I65").Value = 5777.5
$ws.Range("K65").Value = 28887.5
$ws.Range("M65").Value = -25767.5
$ws.Range("H113").Value = 699.2
$ws.Range("I113").Value = 573.4286
$ws.Range("J113").Value = 992.6667
$ws.Range("K113").Value = 573.4286
$ws.Range("L113").Value = 992.6667
$ws.Range("M113").Value = 1596.5714
$ws.Range("N113").Value = -5332.6667
$ws.Range("H122").Value = 2692.3635
$ws.Range("I122").Value = 1670.3334
$ws.Range("J122").Value = 3918.8
$ws.Range("K122").Value = 5011.0002
$ws.Range("L122").Value = 11756.4
$ws.Range("M122").Value = -2561.0002
$ws.Range("N122").Value = -16656.4
$ws.Range("H132").Value = 1970.9656
$ws.Range("I132").Value = 1954.1923
$ws.Range("K132").Value = 5862.5769
$ws.Range("M132").Value = -3332.5769
$ws.Range("H134").Value = 2804.75
$ws.Range("I134").Value = 487.33334
$ws.Range("K134").Value = 1462.00002
$ws.Range("M134").Value = 1072.99998
$ws.Range("H136").Value = 3522.077
$ws.Range("I136").Value = 2871.5908
$ws.Range("K136").Value = 8614.7724
$ws.Range("M136").Value = -6064.7724
$ws.Range("H141").Value = 58016
$ws.Range("J141").Value = 58016
$ws.Range("L141").Value = 58016
$ws.Range("N141").Value = -68376

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 942
$ws.Range("I5").Value = 790.8
$ws.Range("K5").Value = 2372.4
$ws.Range("M5").Value = -2260.4
$ws.Range("H12").Value = 88
$ws.Range("J12").Value = 151.28572
$ws.Range("L12").Value = 453.85716
$ws.Range("N12").Value = -799.85716
$ws.Range("H50").Value = 126747.375
$ws.Range("J50").Value = 799.5
$ws.Range("L50").Value = 2398.5
$ws.Range("N50").Value = -3360.5
$ws.Range("H53").Value = 126747.375
$ws.Range("J53").Value = 799.5
$ws.Range("L53").Value = 2398.5
$ws.Range("N53").Value = -3360.5
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H82").Value = 8925
$ws.Range("J82").Value = 8900
$ws.Range("L82").Value = 26700
$ws.Range("N82").Value = -27512
$ws.Range("H85").Value = 8925
$ws.Range("J85").Value = 8900
$ws.Range("L85").Value = 26700
$ws.Range("N85").Value = -29508
$ws.Range("H106").Value = 20000
$ws.Range("J106").Value = 20000
$ws.Range("L106").Value = 60000
$ws.Range("N106").Value = -61892
$ws.Range("H107").Value = 1098
$ws.Range("J107").Value = 1159.4546
$ws.Range("L107").Value = 3478.3638
$ws.Range("N107").Value = -7318.3638
$ws.Range("H113").Value = 1596.3125
$ws.Range("J113").Value = 1632
$ws.Range("L113").Value = 4896
$ws.Range("N113").Value = -9236
$ws.Range("H131").Value = 13075159
$ws.Range("J131").Value = 10419978
$ws.Range("L131").Value = 31259934
$ws.Range("N131").Value = -31270014
$ws.Range("H132").Value = 9899.75
$ws.Range("J132").Value = 9866.333000000001
$ws.Range("L132").Value = 88796.997
$ws.Range("N132").Value = -93856.997
$ws.Range("H134").Value = 12454.926
$ws.Range("I134").Value = 5254.273
$ws.Range("K134").Value = 15762.819
$ws.Range("M134").Value = -10692.819
$ws.Range("H135").Value = 942
$ws.Range("I135").Value = 790.8
$ws.Range("K135").Value = 7117.2
$ws.Range("M135").Value = -4582.2
$ws.Range("H137").Value = 10717
$ws.Range("J137").Value = 12616.333
$ws.Range("L137").Value = 37848.999
$ws.Range("N137").Value = -48048.999
$ws.Range("H140").Value = 7821173
$ws.Range("I140").Value = 20837474
$ws.Range("J140").Value = 11391.55
$ws.Range("K140").Value = 62512422
$ws.Range("L140").Value = 34174.64999999999
$ws.Range("M140").Value = -62507242
$ws.Range("N140").Value = -44534.64999999999
$ws.Range("H141").Value = 36535.832
$ws.Range("I141").Value = 8772
$ws.Range("K141").Value = 26316
$ws.Range("M141").Value = -21136

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 134.82353
$ws.Range("I2").Value = 50.384617
$ws.Range("J2").Value = 409.25
$ws.Range("K2").Value = 50.384617
$ws.Range("L2").Value = 409.25
$ws.Range("M2").Value = 62.615383
$ws.Range("N2").Value = -635.25
$ws.Range("H15").Value = 59993.5
$ws.Range("J15").Value = 59993.5
$ws.Range("L15").Value = 59993.5
$ws.Range("N15").Value = -60569.5
$ws.Range("H70").Value = 8078.7617
$ws.Range("I70").Value = 7832.5713
$ws.Range("J70").Value = 8571.143
$ws.Range("K70").Value = 7832.5713
$ws.Range("L70").Value = 8571.143
$ws.Range("M70").Value = -7562.5713
$ws.Range("N70").Value = -9111.143
$ws.Range("H73").Value = 8078.7617
$ws.Range("I73").Value = 7832.5713
$ws.Range("J73").Value = 8571.143
$ws.Range("K73").Value = 7832.5713
$ws.Range("L73").Value = 8571.143
$ws.Range("M73").Value = -6896.5713
$ws.Range("N73").Value = -10443.143
$ws.Range("H81").Value = 59993.5
$ws.Range("J81").Value = 59993.5
$ws.Range("L81").Value = 59993.5
$ws.Range("N81").Value = -61989.5
$ws.Range("H84").Value = 59993.5
$ws.Range("J84").Value = 59993.5
$ws.Range("L84").Value = 179980.5
$ws.Range("N84").Value = -189964.5
$ws.Range("H102").Value = 3856
$ws.Range("I102").Value = 3611
$ws.Range("J102").Value = 4284.75
$ws.Range("K102").Value = 3611
$ws.Range("L102").Value = 4284.75
$ws.Range("M102").Value = -1989
$ws.Range("N102").Value = -7528.75
$ws.Range("H122").Value = 1892.2
$ws.Range("I122").Value = 1944.4546
$ws.Range("J122").Value = 1748.5
$ws.Range("K122").Value = 5833.3638
$ws.Range("L122").Value = 5245.5
$ws.Range("M122").Value = -3383.3638
$ws.Range("N122").Value = -10145.5
$ws.Range("H132").Value = 2646.5
$ws.Range("I132").Value = 2992.875
$ws.Range("J132").Value = 2184.6667
$ws.Range("K132").Value = 8978.625
$ws.Range("L132").Value = 6554.000100000001
$ws.Range("M132").Value = -6448.625
$ws.Range("N132").Value = -11614.0001
$ws.Range("H139").Value = 79999
$ws.Range("J139").Value = 79999
$ws.Range("L139").Value = 79999
$ws.Range("N139").Value = -90279

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H8").Value = 265083.5
$ws.Range("J8").Value = 265083.5
$ws.Range("L8").Value = 265083.5
$ws.Range("N8").Value = -265363.5
$ws.Range("H16").Value = 3491
$ws.Range("I16").Value = 4455
$ws.Range("K16").Value = 4455
$ws.Range("M16").Value = -4285
$ws.Range("H46").Value = 2411.2307
$ws.Range("I46").Value = 1700.2858
$ws.Range("J46").Value = 3240.6667
$ws.Range("K46").Value = 1700.2858
$ws.Range("L46").Value = 3240.6667
$ws.Range("M46").Value = -1512.2858
$ws.Range("N46").Value = -3616.6667
$ws.Range("H61").Value = 2233.4285
$ws.Range("I61").Value = 2189.923
$ws.Range("K61").Value = 2189.923
$ws.Range("M61").Value = -1987.923
$ws.Range("H82").Value = 1656.6
$ws.Range("I82").Value = 943.55554
$ws.Range("J82").Value = 2726.1667
$ws.Range("K82").Value = 943.55554
$ws.Range("L82").Value = 2726.1667
$ws.Range("M82").Value = -582.55554
$ws.Range("N82").Value = -3448.1667
$ws.Range("H85").Value = 1656.6
$ws.Range("I85").Value = 943.55554
$ws.Range("J85").Value = 2726.1667
$ws.Range("K85").Value = 943.55554
$ws.Range("L85").Value = 2726.1667
$ws.Range("M85").Value = 304.44446
$ws.Range("N85").Value = -5222.1667
$ws.Range("H87").Value = 98995
$ws.Range("J87").Value = 98995
$ws.Range("L87").Value = 98995
$ws.Range("N87").Value = -101241
$ws.Range("H90").Value = 98995
$ws.Range("J90").Value = 98995
$ws.Range("L90").Value = 296985
$ws.Range("N90").Value = -308217
$ws.Range("H92").Value = 99389
$ws.Range("J92").Value = 99389
$ws.Range("L92").Value = 99389
$ws.Range("N92").Value = -104381
$ws.Range("H100").Value = 1946.0588
$ws.Range("I100").Value = 923.5833
$ws.Range("K100").Value = 923.5833
$ws.Range("M100").Value = -382.5833
$ws.Range("H108").Value = 91001
$ws.Range("J108").Value = 91001
$ws.Range("L108").Value = 91001
$ws.Range("N108").Value = -98681
$ws.Range("H113").Value = 2233.4285
$ws.Range("I113").Value = 2189.923
$ws.Range("K113").Value = 2189.923
$ws.Range("M113").Value = -19.92299999999977
$ws.Range("H118").Value = 100000
$ws.Range("J118").Value = 100000
$ws.Range("L118").Value = 100000
$ws.Range("N118").Value = -103314
$ws.Range("H132").Value = 2431.0977
$ws.Range("I132").Value = 1710.64
$ws.Range("J132").Value = 3556.8125
$ws.Range("K132").Value = 5131.92
$ws.Range("L132").Value = 10670.4375
$ws.Range("M132").Value = -2601.92
$ws.Range("N132").Value = -15730.4375
$ws.Range("H136").Value = 1863.3448
$ws.Range("I136").Value = 1273.1765
$ws.Range("J136").Value = 2699.4167
$ws.Range("K136").Value = 3819.5295
$ws.Range("L136").Value = 8098.250100000001
$ws.Range("M136").Value = -1269.5295
$ws.Range("N136").Value = -13198.2501

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = 0
$ws.Range("H20").Value = 21499
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 21499
$ws.Range("K20").Value = 0
$ws.Range("L20").ClearContents()
$ws.Range("M20").Value = 21499
$ws.Range("N20").Value = -21979
$ws.Range("H23").Value = 15798.8
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 15798.8
$ws.Range("K23").Value = 0
$ws.Range("L23").ClearContents()
$ws.Range("M23").Value = 15798.8
$ws.Range("N23").Value = -16256.8
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").ClearContents()
$ws.Range("N25").Value = 0
$ws.Range("H30").Value = 8011
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()
$ws.Range("H32").Value = 29500
$ws.Range("J32").Value = 29500
$ws.Range("L32").Value = 29500
$ws.Range("N32").Value = -30134
$ws.Range("H62").Value = 8871.143
$ws.Range("I62").Value = 8075
$ws.Range("J62").Value = 9932.666999999999
$ws.Range("K62").Value = 8075
$ws.Range("L62").Value = 9932.666999999999
$ws.Range("M62").Value = -7451
$ws.Range("N62").Value = -11180.667
$ws.Range("H65").Value = 8871.143
$ws.Range("I65").Value = 8075
$ws.Range("J65").Value = 9932.666999999999
$ws.Range("K65").Value = 40375
$ws.Range("L65").Value = 49663.335
$ws.Range("M65").Value = -37255
$ws.Range("N65").Value = -55903.335
$ws.Range("H107").Value = 1133.4688
$ws.Range("I107").Value = 1122.4482
$ws.Range("J107").Value = 1240
$ws.Range("K107").Value = 3367.3446
$ws.Range("L107").Value = 3720
$ws.Range("M107").Value = -1447.3446
$ws.Range("N107").Value = -7560
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").ClearContents()
$ws.Range("N116").Value = 0
$ws.Range("H132").Value = 4148.1113
$ws.Range("I132").Value = 4148.1113
$ws.Range("K132").Value = 12444.3339
$ws.Range("M132").Value = -9914.333899999998
$ws.Range("H136").Value = 5427.8887
$ws.Range("I136").Value = 5570.706
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 16712.118
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -14162.118
$ws.Range("N136").Value = -14100

Write-Output "Applied all edits"